$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(6, 2).Value = "8e8f4ea7d0efe7065f4aa179c66a5993"
$ws.Cells.Item(37, 2).Value = "8cae6f5b968301ba77f87fba8b5b5c6d"
$ws.Cells.Item(53, 2).Value = "b367758740bae7360eac1bc8e5e38bc2"
$ws.Cells.Item(54, 2).Value = "bf44ec96ab80cb1716583fd2713be6b3"
$ws.Cells.Item(58, 2).Value = "e021118948136fc1197f1b99869af114"
$ws.Cells.Item(69, 2).Value = "b0bfef82cb938c7b37b82ac03dd92e02"
$ws.Cells.Item(71, 2).Value = "7d5b1d4c9d76911c7f0629c2bbc3b559"
$ws.Cells.Item(78, 2).Value = "12157f38a2a14ab4bf2cf253cea55772"
$ws.Cells.Item(96, 2).Value = "163a2c95fdc0133f3182e4a2f5981be1"
$ws.Cells.Item(100, 2).Value = "8c38adc983273412ce5a788bae1cd763"
$ws.Cells.Item(107, 2).Value = "c90124b7b564c8fd04454539d0804182"
$ws.Cells.Item(108, 2).Value = "1c6b965ee60990ab7717581b4a83445d"
$ws.Cells.Item(120, 2).Value = "8cef06adee08acc58c2564ba45a92776"
$ws.Cells.Item(130, 2).Value = "cd96b58e7ba840c9698dfaad67319aad"
$ws.Cells.Item(132, 2).Value = "8c1bc713711942151cc4e996d53b960a"
$ws.Cells.Item(142, 2).Value = "575e393b45ab9db58ab117dfedf0f70d"
$ws.Cells.Item(159, 2).Value = "17e6f09fd8ea8a8972bc475df817080f"
$ws.Cells.Item(169, 2).Value = "6afcb86346c0f16cac73003425cae14d"
$ws.Cells.Item(171, 2).Value = "9419ac7b14b927a35392df1206b662a5"
$ws.Cells.Item(173, 2).Value = "c4c5187e346d28891f3aa78ac4ae5d86"
$ws.Cells.Item(190, 2).Value = "a0e66fbb3a80f46243aa89c973e6aef5"
$ws.Cells.Item(233, 2).Value = "4dc6992645510e489bbe6c13b9760931"
$ws.Cells.Item(246, 2).Value = "a7844963b70be534ed450364d9f7d1e9"
$ws.Cells.Item(255, 2).Value = "3c39cc40a5d3c996803a1bbb7835e95b"
$ws.Cells.Item(276, 2).Value = "a5a8399642eb3856bc0ed3d26c605c8e"
$ws.Cells.Item(281, 2).Value = "91d6cecafdef3ad37838abc58fd1f3c8"
$ws.Cells.Item(299, 2).Value = "27b7354351f85b3ec9741b3dc249118a"
$ws.Cells.Item(310, 2).Value = "8154777e2c8ce05773d7088ed02de109"
$ws.Cells.Item(339, 2).Value = "1e506b1f2a033ed20095cbdd53afc20a"
$ws.Cells.Item(343, 2).Value = "9c8e173b79f48d63f00af95644862e76"
$ws.Cells.Item(352, 2).Value = "444d7c36df66c5ffb38e38d0022965ff"
$ws.Cells.Item(361, 2).Value = "a70c9b6abacf1267334eda2a60786805"
$ws.Cells.Item(373, 2).Value = "e48707233d422bf58637c564a378b383"
$ws.Cells.Item(378, 2).Value = "bb5530e3cc2fd9b2c27f6435d21c8d7a"
$ws.Cells.Item(388, 2).Value = "e021a1af0e663045acb12bbf52548523"
$ws.Cells.Item(400, 2).Value = "03a230b2bb153353c297430e4a97c1c0"
$ws.Cells.Item(407, 2).Value = "1eeacbd7d37f53f89db299ee668fff75"
$ws.Cells.Item(410, 2).Value = "c8e5ee6496752aa5375e643c806f31c9"
$ws.Cells.Item(417, 2).Value = "2fe5f54c0a0d39a3106d13918c7f78a7"
$ws.Cells.Item(419, 2).Value = "afba4ee92bb44bede48ddf483ac24705"
$ws.Cells.Item(446, 2).Value = "9de5a67740a3686774a6f39010a19265"
$ws.Cells.Item(452, 2).Value = "81ecf120c44fe8e8d2e0be038b23c315"
$ws.Cells.Item(460, 2).Value = "0cd8625297c32aba25b0f61545f1b53e"
$ws.Cells.Item(471, 2).Value = "7cca2f14b69369550c785fd101af3490"
$ws.Cells.Item(472, 2).Value = "846627bbd541c1508403cdd22739c10b"
$ws.Cells.Item(479, 2).Value = "9c97e798b02676a3ca40e0f0ef22b628"
$ws.Cells.Item(483, 2).Value = "982733bb9ab2f264df7e5a6266d301f8"
$ws.Cells.Item(488, 2).Value = "e0e5a8781dbac31946c52f46dcd95db7"
$ws.Cells.Item(492, 2).Value = "ce84a2a5da4ea27b98021964a91beaa4"
$ws.Cells.Item(500, 2).Value = "59328d6fbee2ac587678815c09af1874"
$ws.Cells.Item(517, 2).Value = "4411e56c2ff7e6ec8787d8f6be166e8b"
$ws.Cells.Item(530, 2).Value = "75bf6026f367a6a3e5c8ad3ab0df4e73"
$ws.Cells.Item(542, 2).Value = "32137b737f73b05333e215fb77c16587"
$ws.Cells.Item(543, 2).Value = "ece6eb734faed0dd6d9b51a279f5053d"
$ws.Cells.Item(561, 2).Value = "c7bc39acd047929c20f71caa2141a1f2"
$ws.Cells.Item(562, 2).Value = "b52bbef24892753b20d5bb0e23e4cee0"
$ws.Cells.Item(566, 2).Value = "93cf8370596863b200b01bd187da9d14"
$ws.Cells.Item(567, 2).Value = "8200d185a168f65afd35873facdc969f"
$ws.Cells.Item(572, 2).Value = "2829c5fc1f67e224165dc8d654e289f4"
$ws.Cells.Item(576, 2).Value = "427e54ec9d1c3da5137eb87a5c650425"
$ws.Cells.Item(588, 2).Value = "70e68d04fe9958f7df543c82e254fe1b"
$ws.Cells.Item(593, 2).Value = "9b9367d22346d83cef61f20fb8cf1f46"
$ws.Cells.Item(627, 2).Value = "cd0f810a0814b71df06adc86d49f9165"
$ws.Cells.Item(628, 2).Value = "b4c28e9a6e235253beea9f6a35999b21"
$ws.Cells.Item(673, 2).Value = "69a991035bad44db6eb9fad5c1f35ab3"
$ws.Cells.Item(686, 2).Value = "d90a7dff71d15179e07e07d8d49cfe8e"
$ws.Cells.Item(727, 2).Value = "ea5085503eeecda17862f1fcddac8e01"
$ws.Cells.Item(734, 2).Value = "2823d56cb3d88595dab1a97de0742c41"
$ws.Cells.Item(749, 2).Value = "b66e0f1588a1a9f9fb6a11b18b4e9d96"
$ws.Cells.Item(756, 2).Value = "9397a483900340432a332a438b43feee"
$ws.Cells.Item(761, 2).Value = "9986aac1f2a947465545084339a92eed"
$ws.Cells.Item(767, 2).Value = "bec68725ca3ed1d2d22a539f7a43ba56"
$ws.Cells.Item(776, 2).Value = "e867a7ef5a2c4abfe453536e28f5ee67"
$ws.Cells.Item(786, 2).Value = "d7c32f6feaa74b68ad82f3fb3036d04e"
$ws.Cells.Item(819, 2).Value = "f918429f8f38492013789bfd11f54108"
$ws.Cells.Item(824, 2).Value = "acf46416e7e91b24f3fa645c04779926"
$ws.Cells.Item(848, 2).Value = "ae7efecaf8736ca69f95c36d2f77d0d1"
$ws.Cells.Item(855, 2).Value = "ec5110340224ff40e879ea2857e85751"
$ws.Cells.Item(856, 2).Value = "ba87f07a3bf7bdee814ce70142401d55"
$ws.Cells.Item(869, 2).Value = "87d5f4401301379682bc0ad75b7a1ef8"
$ws.Cells.Item(881, 2).Value = "57e6a5c84dd98419a196f687ba7da1cf"
$ws.Cells.Item(902, 2).Value = "b98caff4d64e24b295aa661bef3da148"
$ws.Cells.Item(904, 2).Value = "cd1a090fd82a983cf3eef5f74f74fdd1"
$ws.Cells.Item(928, 2).Value = "075dc0b3177c298bc5836ccf2890df11"
$ws.Cells.Item(938, 2).Value = "cddf31685c87ea1f10494c5af2604c98"
$ws.Cells.Item(952, 2).Value = "56434335b1d337bea2b9eba3b059519e"
$ws.Cells.Item(953, 2).Value = "4f5e17e055f48fc2357151abfc4241f0"
